{"js": "// Replace the date line and every \"A\u00d7B=C\" equation in the table with the\n// updated values from the commit. Each \"before\" value is unique in the\n// document, so a scoped search-and-replace on context.document.body is\n// unambiguous and safe.\nconst replacements = [\n  [\"2024-02-10 Saturday\", \"2024-02-11 Sunday\"],\n  [\"670\u00d76=4020\", \"392\u00d79=3528\"],\n  [\"394\u00d77=2758\", \"174\u00d75=870\"],\n  [\"268\u00d73=804\", \"424\u00d73=1272\"],\n  [\"704\u00d76=4224\", \"802\u00d76=4812\"],\n  [\"468\u00d75=2340\", \"514\u00d76=3084\"],\n  [\"989\u00d77=6923\", \"904\u00d72=1808\"],\n  [\"579\u00d79=5211\", \"518\u00d79=4662\"],\n  [\"946\u00d75=4730\", \"139\u00d75=695\"],\n  [\"882\u00d75=4410\", \"431\u00d72=862\"],\n  [\"652\u00d73=1956\", \"299\u00d72=598\"],\n  [\"147\u00d76=882\", \"333\u00d74=1332\"],\n  [\"939\u00d78=7512\", \"422\u00d72=844\"],\n  [\"680\u00d72=1360\", \"115\u00d79=1035\"],\n  [\"284\u00d76=1704\", \"338\u00d76=2028\"],\n  [\"176\u00d73=528\", \"466\u00d73=1398\"],\n  [\"441\u00d76=2646\", \"985\u00d75=4925\"],\n  [\"675\u00d76=4050\", \"975\u00d75=4875\"],\n  [\"346\u00d76=2076\", \"298\u00d73=894\"],\n  [\"990\u00d78=7920\", \"238\u00d79=2142\"],\n  [\"334\u00d76=2004\", \"506\u00d75=2530\"],\n  [\"931\u00d72=1862\", \"701\u00d78=5608\"],\n  [\"772\u00d77=5404\", \"204\u00d79=1836\"],\n  [\"422\u00d78=3376\", \"437\u00d79=3933\"],\n  [\"828\u00d76=4968\", \"861\u00d74=3444\"],\n  [\"652\u00d74=2608\", \"985\u00d72=1970\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${before}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"A\u00d7B=C\" equation in the table with the\n# updated values from the commit. Each \"before\" value is unique in the\n# document, so Find/Replace scoped to the whole document body is\n# unambiguous and safe.\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$pairs = @(\n    @(\"2024-02-10 Saturday\", \"2024-02-11 Sunday\"),\n    @(\"670\u00d76=4020\", \"392\u00d79=3528\"),\n    @(\"394\u00d77=2758\", \"174\u00d75=870\"),\n    @(\"268\u00d73=804\", \"424\u00d73=1272\"),\n    @(\"704\u00d76=4224\", \"802\u00d76=4812\"),\n    @(\"468\u00d75=2340\", \"514\u00d76=3084\"),\n    @(\"989\u00d77=6923\", \"904\u00d72=1808\"),\n    @(\"579\u00d79=5211\", \"518\u00d79=4662\"),\n    @(\"946\u00d75=4730\", \"139\u00d75=695\"),\n    @(\"882\u00d75=4410\", \"431\u00d72=862\"),\n    @(\"652\u00d73=1956\", \"299\u00d72=598\"),\n    @(\"147\u00d76=882\", \"333\u00d74=1332\"),\n    @(\"939\u00d78=7512\", \"422\u00d72=844\"),\n    @(\"680\u00d72=1360\", \"115\u00d79=1035\"),\n    @(\"284\u00d76=1704\", \"338\u00d76=2028\"),\n    @(\"176\u00d73=528\", \"466\u00d73=1398\"),\n    @(\"441\u00d76=2646\", \"985\u00d75=4925\"),\n    @(\"675\u00d76=4050\", \"975\u00d75=4875\"),\n    @(\"346\u00d76=2076\", \"298\u00d73=894\"),\n    @(\"990\u00d78=7920\", \"238\u00d79=2142\"),\n    @(\"334\u00d76=2004\", \"506\u00d75=2530\"),\n    @(\"931\u00d72=1862\", \"701\u00d78=5608\"),\n    @(\"772\u00d77=5404\", \"204\u00d79=1836\"),\n    @(\"422\u00d78=3376\", \"437\u00d79=3933\"),\n    @(\"828\u00d76=4968\", \"861\u00d74=3444\"),\n    @(\"652\u00d74=2608\", \"985\u00d72=1970\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute($findText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\nWrite-Output \"done\"\n"}
